$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.050.08"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "1.831.04"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.40"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6542"
$ws.Range("E6").Value = "  -2.97%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.48"
$ws.Range("E8").Value = "  +5.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2938"
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07344"
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.98"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07674"
$ws.Range("E12").Value = "  -0.57%  "
$ws.Range("D13").Value = "1.834.88"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6674"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.88"
$ws.Range("E16").Value = "  -4.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.104"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("E18").Value = "  +4.55%  "
$ws.Range("D19").Value = "29.053.15"
$ws.Range("E19").Value = "  -0.85%  "
$ws.Range("D20").Value = "2.085.95"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.44"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.25"
$ws.Range("E22").Value = "  -1.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.114"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.23"
$ws.Range("E26").Value = "  -1.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.504"
$ws.Range("E27").Value = "  -2.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1378"
$ws.Range("E28").Value = "  -1.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.89"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.508"
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.110"
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.017"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.199"
$ws.Range("E33").Value = "  +0.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05350"
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7436"
$ws.Range("E35").Value = "  -2.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.840"
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.155"
$ws.Range("E37").Value = "  +1.93%  "
$ws.Range("E38").Value = "  -1.26%  "
$ws.Range("D39").Value = "1.298.67"
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01787"
$ws.Range("E40").Value = "  -0.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.748"
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.362"
$ws.Range("E42").Value = "  +7.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8988"
$ws.Range("E43").Value = "  -2.36%  "
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.23"
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("D46").Value = "1.983.99"
$ws.Range("E46").Value = "  -0.87%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5139"
$ws.Range("E47").Value = "  -0.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.21"
$ws.Range("E48").Value = "  +0.66%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000120"
$ws.Range("E49").Value = "  -5.22%  "
$ws.Range("B50").Value = "XinFinNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07625"
$ws.Range("E50").Value = "  -6.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.742"
$ws.Range("E51").Value = "  -1.90%  "
